$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = $null
$ws.Range("H15").Value = 1592.683
$ws.Range("I15").Value = 1592.683
$ws.Range("K15").Value = 4778.049
$ws.Range("M15").Value = -4609.049
$ws.Range("H19").Value = 767
$ws.Range("I19").Value = 674
$ws.Range("J19").Value = 860
$ws.Range("K19").Value = 674
$ws.Range("L19").Value = 860
$ws.Range("M19").Value = -499
$ws.Range("N19").Value = -1210
$ws.Range("H58").Value = 938.9167
$ws.Range("J58").Value = 1569.5714
$ws.Range("L58").Value = 4708.7142
$ws.Range("N58").Value = -5008.7142
$ws.Range("H64").Value = 5156.25
$ws.Range("I64").Value = 4200
$ws.Range("K64").Value = 4200
$ws.Range("M64").Value = -3952
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 5156.25
$ws.Range("I67").Value = 4200
$ws.Range("K67").Value = 4200
$ws.Range("M67").Value = -3342
$ws.Range("N67").Value = $null
$ws.Range("H69").Value = 7265.675
$ws.Range("I69").Value = 3000
$ws.Range("K69").Value = 9000
$ws.Range("M69").Value = -8126
$ws.Range("H72").Value = 7265.675
$ws.Range("I72").Value = 3000
$ws.Range("K72").Value = 27000
$ws.Range("M72").Value = -22632
$ws.Range("H95").Value = 36662.668
$ws.Range("J95").Value = 36662.668
$ws.Range("L95").Value = 36662.668
$ws.Range("N95").Value = -42154.668
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = 754
$ws.Range("H116").Value = 20687.5
$ws.Range("I116").Value = 20687.5
$ws.Range("K116").Value = 20687.5
$ws.Range("M116").Value = -17245.5
$ws.Range("H132").Value = 17485.295
$ws.Range("I132").Value = 18988.334
$ws.Range("K132").Value = 56965.00199999999
$ws.Range("M132").Value = -54435.00199999999
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3032836.5
$ws.Range("I32").Value = 540.5185
$ws.Range("K32").Value = 540.5185
$ws.Range("M32").Value = -253.5185
$ws.Range("N32").Value = $null
$ws.Range("H45").Value = 2188.3125
$ws.Range("I45").Value = 1737.5555
$ws.Range("J45").Value = 2767.8572
$ws.Range("K45").Value = 1737.5555
$ws.Range("L45").Value = 2767.8572
$ws.Range("M45").Value = -1360.5555
$ws.Range("N45").Value = -3521.8572
$ws.Range("H63").Value = 2020.7142
$ws.Range("I63").Value = 2020.7142
$ws.Range("K63").Value = 2020.7142
$ws.Range("M63").Value = -1334.7142
$ws.Range("H66").Value = 2020.7142
$ws.Range("I66").Value = 2020.7142
$ws.Range("K66").Value = 10103.571
$ws.Range("M66").Value = -6671.571
$ws.Range("H101").Value = 67500
$ws.Range("J101").Value = 67500
$ws.Range("L101").Value = 67500
$ws.Range("N101").Value = -73990
$ws.Range("H124").Value = 28950
$ws.Range("J124").Value = 28950
$ws.Range("L124").Value = 28950
$ws.Range("N124").Value = -38770
$ws.Range("H132").Value = 3989.0833
$ws.Range("I132").Value = 3989.0833
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11967.2499
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9437.249899999999
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 550
$ws.Range("I20").Value = 550
$ws.Range("K20").Value = 550
$ws.Range("M20").Value = -303
$ws.Range("H86").Value = 2826.5264
$ws.Range("I86").Value = 1350.4166
$ws.Range("J86").Value = 5357
$ws.Range("K86").Value = 1350.4166
$ws.Range("L86").Value = 5357
$ws.Range("M86").Value = -227.4166
$ws.Range("N86").Value = -7603
$ws.Range("H89").Value = 2826.5264
$ws.Range("I89").Value = 1350.4166
$ws.Range("J89").Value = 5357
$ws.Range("K89").Value = 6752.083000000001
$ws.Range("L89").Value = 26785
$ws.Range("M89").Value = -1136.083000000001
$ws.Range("N89").Value = -38017
$ws.Range("H99").Value = 76923940
$ws.Range("I99").Value = 90909910
$ws.Range("J99").Value = 1082.5
$ws.Range("K99").Value = 90909910
$ws.Range("L99").Value = 1082.5
$ws.Range("M99").Value = -90908412
$ws.Range("N99").Value = -4078.5
$ws.Range("H105").Value = 6061941.5
$ws.Range("I105").Value = 9092052
$ws.Range("K105").Value = 9092052
$ws.Range("M105").Value = -9090305
$ws.Range("N105").Value = $null
$ws.Range("H134").Value = 3803
$ws.Range("I134").Value = 1156.1428
$ws.Range("K134").Value = 3468.4284
$ws.Range("M134").Value = -933.4284000000002
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3704.5625
$ws.Range("J31").Value = 4631.579
$ws.Range("L31").Value = 4631.579
$ws.Range("N31").Value = -5221.579
$ws.Range("H34").Value = 3704.5625
$ws.Range("J34").Value = 4631.579
$ws.Range("L34").Value = 4631.579
$ws.Range("N34").Value = -5035.579
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H132").Value = 2356.44
$ws.Range("I132").Value = 2329.625
$ws.Range("K132").Value = 6988.875
$ws.Range("M132").Value = -4458.875
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22.324324
$ws.Range("J2").Value = 21.8125
$ws.Range("L2").Value = 130.875
$ws.Range("N2").Value = -356.875
$ws.Range("H4").Value = 2567923
$ws.Range("I4").Value = 2567923
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 7703769
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -7703657
$ws.Range("N4").Value = $null
$ws.Range("H130").Value = 2663.3333
$ws.Range("I130").Value = 2663.3333
$ws.Range("K130").Value = 7989.999899999999
$ws.Range("M130").Value = -2969.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 5000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null
$ws.Range("H102").Value = 1861.64
$ws.Range("I102").Value = 1833.3182
$ws.Range("J102").Value = 2069.3333
$ws.Range("K102").Value = 1833.3182
$ws.Range("L102").Value = 2069.3333
$ws.Range("M102").Value = -211.3181999999999
$ws.Range("N102").Value = -5313.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2217.842
$ws.Range("I46").Value = 590.7143
$ws.Range("J46").Value = 3167
$ws.Range("K46").Value = 590.7143
$ws.Range("L46").Value = 3167
$ws.Range("M46").Value = -402.7143
$ws.Range("N46").Value = -3543
$ws.Range("H82").Value = 3135.6428
$ws.Range("I82").Value = 1143.75
$ws.Range("K82").Value = 1143.75
$ws.Range("M82").Value = -782.75
$ws.Range("N82").Value = $null
$ws.Range("H85").Value = 3135.6428
$ws.Range("I85").Value = 1143.75
$ws.Range("K85").Value = 1143.75
$ws.Range("M85").Value = 104.25
$ws.Range("N85").Value = $null
$ws.Range("H136").Value = 1975.1428
$ws.Range("I136").Value = 1679.3334
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 5038.0002
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -2488.0002
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = $null
$ws.Range("N38").Value = $null
$ws.Range("H43").Value = 52500
$ws.Range("I43").Value = 45000
$ws.Range("J43").Value = 60000
$ws.Range("K43").Value = 45000
$ws.Range("L43").Value = 60000
$ws.Range("M43").Value = -44851
$ws.Range("N43").Value = -60298
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = $null
$ws.Range("H122").Value = 2078.5625
$ws.Range("I122").Value = 1246.7273
$ws.Range("J122").Value = 3908.6
$ws.Range("K122").Value = 3740.1819
$ws.Range("L122").Value = 11725.8
$ws.Range("M122").Value = -1290.1819
$ws.Range("N122").Value = -16625.8
$ws.Range("H132").Value = 2160.2307
$ws.Range("I132").Value = 1808.4
$ws.Range("J132").Value = 3333
$ws.Range("K132").Value = 5425.200000000001
$ws.Range("L132").Value = 9999
$ws.Range("M132").Value = -2895.200000000001
$ws.Range("N132").Value = -15059
